# Daily attendance update - 2025-08-25
#
# Column AD on the "WCS_Team_August_2025" sheet is the date column for
# Monday, 2025-08-25. Fill in each employee's attendance status for that
# day (rows 3-18), re-using the "Good" (WFO) / "Neutral" (SL) cell
# formatting already applied earlier in each employee's row (column I
# carries the same status for that employee, already styled correctly).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WCS_Team_August_2025")

# Employee -> status for 2025-08-25 (column AD), rows 3-18
$statuses = @{
    3  = "WFO"
    4  = "WFO"
    5  = "WFO"
    6  = "WFO"
    7  = "WFO"
    8  = "WFO"
    9  = "WFO"
    10 = "WFO"
    11 = "WFO"
    12 = "WFO"
    13 = "WFO"
    14 = "WFO"
    15 = "WFO"
    16 = "SL"
    17 = "WFO"
    18 = "WFO"
}

foreach ($row in 3..18) {
    $status = $statuses[$row]

    # Set the value first (so dependent COUNTIF formulas recalc).
    $dst = $ws.Range("AD$row")
    $dst.Value = $status

    # Column I of the same row already has an identically-valued cell
    # with the correct conditional "Good"/"Neutral" formatting - copy
    # just its formatting (not its value) onto the new cell.
    $src = $ws.Range("I$row")
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# Leave the sheet selecting the range that was just filled in.
$ws.Activate() | Out-Null
$ws.Range("AD3:AD18").Select() | Out-Null
